$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Top 10 2017 -> 2021 mapping text (content refresh: EN->ES fixes,
# reworded category names, "(Nuevo)" -> "(Nueva)") ---

$ws.Range("C7").Value = "A01:2017-Injection"
$ws.Range("C8").Value = "A02:2017-Broken Authentication"
$ws.Range("C9").Value = "A03:2017-Sensitive Data Exposure"
$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"

$ws.Range("D10").Value = "(Nueva)"
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("D16").Value = "(Nueva)"

$ws.Range("E8").Value = "A02:2021-Fallas Criptográficas"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"

# --- Remove the two stray duplicate connector shapes left over in the drawing ---
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()

# --- Column sizing: drop helper column B's explicit width, narrow column C ---
$ws.Columns.Item(3).ColumnWidth = 47.666666666666664

# --- Selection moved to E23 ---
$ws.Range("E23").Select()
